$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "FailedPostOutput" column (E) entirely, shifting remaining
# columns left.
$ws.Range("E1").EntireColumn.Delete()

# Rename "NetworkConnectivity" -> "NetworkConnectivitySNV-US1".
# After the deletion above this header moved from J1 to I1.
$ws.Range("I1").Value = "NetworkConnectivitySNV-US1"

# Remove the trailing "TimeSinceSystemIsUp" column, which is now the last
# used column (K).
$ws.Range("K1").EntireColumn.Delete()
